$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Homework column (H): remove "Homework Reflection N" entries ---
$ws.Range("H9").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("H23").ClearContents()
$ws.Range("H31").ClearContents()

# --- Project column (I): remove proposal-related milestones ---
$ws.Range("I11").ClearContents()
$ws.Range("I14").ClearContents()

# --- Project column (I): shift remaining milestones up two rows ---
$ws.Range("I13").Value = $ws.Range("I17").Value2
$ws.Range("I17").ClearContents()

$ws.Range("I15").Value = $ws.Range("I20").Value2
$ws.Range("I20").ClearContents()

$ws.Range("I22").Value = $ws.Range("I23").Value2
$ws.Range("I23").ClearContents()

$ws.Range("I24").Value = $ws.Range("I26").Value2
$ws.Range("I26").ClearContents()

# --- View: reset scroll position and selection ---
$ws.Range("D3").Select()
